# Modelo R.xlsx — "Se añadieron los pagos y se ha hecho una pequeña
# modificacion en la base de datos."
#
# Schema-diagram sheet: the "Depto" table loses its IDContrato column
# (that relation now lives on "pago", where it already existed), and the
# "Contrato" table gains a new IDDepto column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Depto table (row 5): drop the IDContrato column (F5). Use Clear() (not
# ClearContents) so the cell node itself disappears and the row's used
# range shrinks from 2:6 down to 2:5, matching the rest of the Depto block.
$ws.Range("F5").Clear()

# Contrato table (row 21): add a new IDDepto column in E21. Copy the
# formatting from the neighboring D21 cell first (so it picks up the
# same thin-border style), then overwrite the value with the new label.
$ws.Range("D21").Copy($ws.Range("E21"))
$ws.Range("E21").Value2 = "IDDepto"

# Move the active selection to match the saved cursor position in the
# workbook (was G15, now H7).
$ws.Range("H7").Select()
